$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows into the condition table ---
# Duplicate row 4 (which carries the "s=1" cell style used by the labelled
# rows) and insert the copy above row 3. This creates a new row 3 and
# pushes the former rows 3 and 4 down to rows 4 and 5.
$ws.Rows.Item(4).Copy()
$ws.Rows.Item(3).Insert()

# Duplicate row 5 (now a copy of the old "blank" row) and insert the copy
# above row 5 itself, giving a new row 5 and pushing the former row 5 (the
# original "blank" row) down to row 6.
$ws.Rows.Item(5).Copy()
$ws.Rows.Item(5).Insert()

# --- Row 2: cont1_hi -- refreshed trial parameters ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "cont1_hi"
$ws.Range("C2").Value = 0.1
$ws.Range("D2").Value = 0.02
$ws.Range("E2").Value = 0.01
$ws.Range("F2").Value = "colour"
$ws.Range("G2").Value = 0.1
$ws.Range("H2").Value = 500
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 0.5
$ws.Range("K2").Value = 0.2
$ws.Range("L2").Value = 12
$ws.Range("M2").Value = 250
$ws.Range("N2").Value = 10
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 17

# --- Row 3 (new): cont1_lo ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "cont1_lo"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = 0.01
$ws.Range("D3").Value = 0.02
$ws.Range("E3").Value = 0.01
$ws.Range("F3").Value = "colour"
$ws.Range("G3").Value = 0.1
$ws.Range("H3").Value = 500
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 0.5
$ws.Range("K3").Value = 0.2
$ws.Range("L3").Value = 12
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 10
$ws.Range("O3").Value = 2
$ws.Range("P3").Value = 17

# --- Row 4: cont7_hi -- refreshed trial parameters ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "cont7_hi"
$ws.Range("C4").Value = 0.1
$ws.Range("D4").Value = 0.02
$ws.Range("E4").Value = 0.01
$ws.Range("F4").Value = "colour"
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 0.5
$ws.Range("K4").Value = 0.2
$ws.Range("L4").Value = 12
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 10
$ws.Range("O4").Value = 2
$ws.Range("P4").Value = 17

# --- Row 5 (new): cont7_lo ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "cont7_lo"
$ws.Range("C5").Value = 0.01
$ws.Range("D5").Value = 0.02
$ws.Range("E5").Value = 0.01
$ws.Range("F5").Value = "colour"
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 500
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 0.5
$ws.Range("K5").Value = 0.2
$ws.Range("L5").Value = 12
$ws.Range("M5").Value = 250
$ws.Range("N5").Value = 10
$ws.Range("O5").Value = 2
$ws.Range("P5").Value = 17

# --- Row 6: blank -- refreshed trial parameters ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "blank"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0.001
$ws.Range("E6").Value = 0.001
$ws.Range("F6").Value = "colour"
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 500
$ws.Range("I6").Value = 1000
$ws.Range("J6").Value = 0.5
$ws.Range("K6").Value = 0.2
$ws.Range("L6").Value = 12
$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 10
$ws.Range("O6").Value = 2
$ws.Range("P6").Value = 17

# --- Sheet view bookkeeping ---
$ws.Range("C5").Select() | Out-Null
